$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2309.4546
$ws.Range("I40").Value = 2114.2856
$ws.Range("J40").Value = 2651
$ws.Range("K40").Value = 2114.2856
$ws.Range("L40").Value = 2651
$ws.Range("M40").Value = -1939.2856
$ws.Range("N40").Value = -3001
$ws.Range("H52").Value = 1569.6666
$ws.Range("I52").Value = 999.5
$ws.Range("J52").Value = 2710
$ws.Range("K52").Value = 2998.5
$ws.Range("L52").Value = 8130
$ws.Range("M52").Value = -2838.5
$ws.Range("N52").Value = -8450
$ws.Range("H58").Value = 1294.7273
$ws.Range("J58").Value = 2803.4
$ws.Range("L58").Value = 8410.200000000001
$ws.Range("N58").Value = -8710.200000000001
$ws.Range("H64").Value = 4896.1333
$ws.Range("I64").Value = 3563.25
$ws.Range("J64").Value = 10227.667
$ws.Range("K64").Value = 3563.25
$ws.Range("L64").Value = 10227.667
$ws.Range("M64").Value = -3315.25
$ws.Range("N64").Value = -10723.667
$ws.Range("H67").Value = 4896.1333
$ws.Range("I67").Value = 3563.25
$ws.Range("J67").Value = 10227.667
$ws.Range("K67").Value = 3563.25
$ws.Range("L67").Value = 10227.667
$ws.Range("M67").Value = -2705.25
$ws.Range("N67").Value = -11943.667
$ws.Range("H70").Value = 3077.6667
$ws.Range("I70").Value = 3000
$ws.Range("J70").Value = 3116.5
$ws.Range("K70").Value = 9000
$ws.Range("L70").Value = 9349.5
$ws.Range("M70").Value = -8730
$ws.Range("N70").Value = -9889.5
$ws.Range("H73").Value = 3077.6667
$ws.Range("I73").Value = 3000
$ws.Range("J73").Value = 3116.5
$ws.Range("K73").Value = 9000
$ws.Range("L73").Value = 9349.5
$ws.Range("M73").Value = -8064
$ws.Range("N73").Value = -11221.5
$ws.Range("H76").Value = 5053559.5
$ws.Range("I76").Value = 5053559.5
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 5053559.5
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -5053244.5
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 5053559.5
$ws.Range("I79").Value = 5053559.5
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 5053559.5
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -5052467.5
$ws.Range("N79").ClearContents()
$ws.Range("H82").Value = 2920.2856
$ws.Range("I82").Value = 2088.4
$ws.Range("J82").Value = 5000
$ws.Range("K82").Value = 6265.200000000001
$ws.Range("L82").Value = 15000
$ws.Range("M82").Value = -5859.200000000001
$ws.Range("N82").Value = -15812
$ws.Range("H85").Value = 2920.2856
$ws.Range("I85").Value = 2088.4
$ws.Range("J85").Value = 5000
$ws.Range("K85").Value = 6265.200000000001
$ws.Range("L85").Value = 15000
$ws.Range("M85").Value = -4861.200000000001
$ws.Range("N85").Value = -17808
$ws.Range("H88").Value = 18521002
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 22224802
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 22224802
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -22225614
$ws.Range("H91").Value = 18521002
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 22224802
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 22224802
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -22227610
$ws.Range("H97").Value = 1000
$ws.Range("J97").Value = 1000
$ws.Range("L97").Value = 3000
$ws.Range("N97").Value = -3992
$ws.Range("H100").Value = 30304750
$ws.Range("I100").Value = 20834334
$ws.Range("J100").Value = 55559190
$ws.Range("K100").Value = 20834334
$ws.Range("L100").Value = 55559190
$ws.Range("M100").Value = -20833793
$ws.Range("N100").Value = -55560272
$ws.Range("H103").Value = 940.25
$ws.Range("I103").Value = 500
$ws.Range("K103").Value = 1500
$ws.Range("M103").Value = -914
$ws.Range("H106").Value = 7475135
$ws.Range("I106").Value = 8009016
$ws.Range("K106").Value = 8009016
$ws.Range("M106").Value = -8008385
$ws.Range("H115").Value = 1328.75
$ws.Range("I115").Value = 1328.75
$ws.Range("K115").Value = 3986.25
$ws.Range("M115").Value = -2419.25
$ws.Range("H118").Value = 463.75
$ws.Range("I118").Value = 387.14285
$ws.Range("J118").Value = 1000
$ws.Range("K118").Value = 1161.42855
$ws.Range("L118").Value = 3000
$ws.Range("M118").Value = 495.5714499999999
$ws.Range("N118").Value = -6314
$ws.Range("H132").Value = 331071.97
$ws.Range("I132").Value = 450375.2
$ws.Range("J132").Value = 38236.816
$ws.Range("K132").Value = 1351125.6
$ws.Range("L132").Value = 114710.448
$ws.Range("M132").Value = -1348595.6
$ws.Range("N132").Value = -119770.448

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 50001
$ws.Range("J6").Value = 50000
$ws.Range("L6").Value = 50000
$ws.Range("N6").Value = -50346
$ws.Range("H9").Value = 19000
$ws.Range("I9").Value = 50000
$ws.Range("J9").Value = 11250
$ws.Range("K9").Value = 50000
$ws.Range("L9").Value = 11250
$ws.Range("M9").Value = -49830
$ws.Range("N9").Value = -11590
$ws.Range("H12").Value = 25022.5
$ws.Range("I12").Value = 25022.5
$ws.Range("K12").Value = 25022.5
$ws.Range("M12").Value = -24849.5
$ws.Range("H20").Value = 19000
$ws.Range("I20").Value = 50000
$ws.Range("J20").Value = 11250
$ws.Range("K20").Value = 50000
$ws.Range("L20").Value = 11250
$ws.Range("M20").Value = -49730
$ws.Range("N20").Value = -11790
$ws.Range("H74").Value = 5051.8184
$ws.Range("J74").Value = 13215
$ws.Range("L74").Value = 13215
$ws.Range("N74").Value = -14963
$ws.Range("H77").Value = 5051.8184
$ws.Range("J77").Value = 13215
$ws.Range("L77").Value = 66075
$ws.Range("N77").Value = -74811
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1231.49
$ws.Range("I68").Value = 966.0909
$ws.Range("J68").Value = 1555.8667
$ws.Range("K68").Value = 2898.2727
$ws.Range("L68").Value = 4667.6001
$ws.Range("M68").Value = -2087.2727
$ws.Range("N68").Value = -6289.6001
$ws.Range("H71").Value = 1231.49
$ws.Range("I71").Value = 966.0909
$ws.Range("J71").Value = 1555.8667
$ws.Range("K71").Value = 8694.8181
$ws.Range("L71").Value = 14002.8003
$ws.Range("M71").Value = -4638.8181
$ws.Range("N71").Value = -22114.8003
$ws.Range("H103").Value = 1682.5238
$ws.Range("I103").Value = 717
$ws.Range("J103").Value = 2406.6667
$ws.Range("K103").Value = 2151
$ws.Range("L103").Value = 7220.000100000001
$ws.Range("M103").Value = -1272
$ws.Range("N103").Value = -8978.000100000001
$ws.Range("H129").Value = 1443.5
$ws.Range("J129").Value = 1899.8889
$ws.Range("L129").Value = 5699.6667
$ws.Range("N129").Value = -15699.6667
$ws.Range("H131").Value = 3416
$ws.Range("J131").Value = 3587.8542
$ws.Range("L131").Value = 10763.5626
$ws.Range("N131").Value = -20843.5626
$ws.Range("H133").Value = 12000
$ws.Range("J133").Value = 14000
$ws.Range("L133").Value = 42000
$ws.Range("N133").Value = -52120
$ws.Range("H134").Value = 49865.715
$ws.Range("I134").Value = 3265
$ws.Range("J134").Value = 112000
$ws.Range("K134").Value = 9795
$ws.Range("L134").Value = 336000
$ws.Range("M134").Value = -4725
$ws.Range("N134").Value = -346140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 30000
$ws.Range("J15").Value = 30000
$ws.Range("L15").Value = 30000
$ws.Range("N15").Value = -30576
$ws.Range("H81").Value = 30000
$ws.Range("J81").Value = 30000
$ws.Range("L81").Value = 30000
$ws.Range("N81").Value = -31996
$ws.Range("H84").Value = 30000
$ws.Range("J84").Value = 30000
$ws.Range("L84").Value = 90000
$ws.Range("N84").Value = -99984

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1030.8889
$ws.Range("I22").Value = 1019.1429
$ws.Range("J22").Value = 1072
$ws.Range("K22").Value = 1019.1429
$ws.Range("L22").Value = 1072
$ws.Range("M22").Value = -724.1429000000001
$ws.Range("N22").Value = -1662
$ws.Range("H27").Value = 1030.8889
$ws.Range("I27").Value = 1019.1429
$ws.Range("J27").Value = 1072
$ws.Range("K27").Value = 1019.1429
$ws.Range("L27").Value = 1072
$ws.Range("M27").Value = -912.1429000000001
$ws.Range("N27").Value = -1286
$ws.Range("H36").Value = 20275.445
$ws.Range("J36").Value = 20275.445
$ws.Range("L36").Value = 20275.445
$ws.Range("N36").Value = -21399.445

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 29375
$ws.Range("J28").Value = 22500
$ws.Range("L28").Value = 22500
$ws.Range("N28").Value = -23196
$ws.Range("H86").Value = 14355.556
$ws.Range("J86").Value = 14355.556
$ws.Range("L86").Value = 14355.556
$ws.Range("N86").Value = -16601.556
$ws.Range("H89").Value = 14355.556
$ws.Range("J89").Value = 14355.556
$ws.Range("L89").Value = 71777.78
$ws.Range("N89").Value = -83009.78
$ws.Range("H136").Value = 8574677
$ws.Range("I136").Value = 23881568
$ws.Range("J136").Value = 2817.92
$ws.Range("K136").Value = 71644704
$ws.Range("L136").Value = 8453.76
$ws.Range("M136").Value = -71642154
$ws.Range("N136").Value = -13553.76

